# Fixed bug with plea check for dismissed.
# Clear the "No Contest" plea and "Guilty" finding cells, and zero out
# the fine amount and suspended fine amount values.

$d = $word.ActiveDocument

$d.Content.Find.Execute("No Contest", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)

$d.Content.Find.Execute("Guilty", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)

$d.Content.Find.Execute("50", $true, $false, $false, $false, $false,
                         $true, 1, $false, "0", 2)

$d.Content.Find.Execute("25", $true, $false, $false, $false, $false,
                         $true, 1, $false, "0", 2)
